$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new EUR->ARS quote row (row 74). The date/time columns look
# like dates/times to Excel's auto-detection, so force them to be written
# as literal text (temporarily marking the cell as Text, then clearing the
# format back off again so the new row doesn't end up with a different
# style than the rest of the sheet).

$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = "2025-10-12"
$ws.Range("A74").ClearFormats()

$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "21:19:37"
$ws.Range("B74").ClearFormats()

$ws.Range("C74").Value = "1.00 EUR = 1,756.2048"
